# Insert a new data row at row 72 (pushing the existing rows 72-75 down to
# 73-76) and populate the new row with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 72 and below down by one row.
$ws.Rows(72).Insert()

# Fill in the new row 72 with the new record.
$ws.Cells.Item(72, 1).Value = 4
$ws.Cells.Item(72, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(72, 3).Value = "Los Lagos"
$ws.Cells.Item(72, 4).Value = 44516
$ws.Cells.Item(72, 5).Value = 10
$ws.Cells.Item(72, 6).Value = 100112052
$ws.Cells.Item(72, 7).Value = "Albahaca"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 100
$ws.Cells.Item(72, 11).Value = 8000
$ws.Cells.Item(72, 12).Value = 8000
$ws.Cells.Item(72, 13).Value = 8000
$ws.Cells.Item(72, 14).Value = "`$/docena de matas"
$ws.Cells.Item(72, 15).Value = "Región Metropolitana"
$ws.Cells.Item(72, 16).Value = 1333
$ws.Cells.Item(72, 17).Value = 6
$ws.Cells.Item(72, 18).Value = "Hortaliza"
